$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q2").Value = 0.7797696841910529
$ws.Range("Q3").Value = 0.7842830816747302
$ws.Range("Q4").Value = 0.7605561072558288
$ws.Range("Q5").Value = 0.7755512153439872
$ws.Range("Q6").Value = 0.8044430831097903
$ws.Range("Q7").Value = 0.7737301978604616
